$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date-formatted style (same as A283:A288) down through the new rows
$ws.Range("A288").Copy()
$ws.Range("A289:A309").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$r = 289
$ws.Cells.Item($r, 1).Value = 45828
$ws.Cells.Item($r, 2).Value = "Flowering"
$ws.Cells.Item($r, 3).Value = "Large"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 85
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 2
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 7
$ws.Cells.Item($r, 13).Value = 0.65
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 29.98
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 0.2
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 50
$ws.Cells.Item($r, 20).Value = 12

$r = 290
$ws.Cells.Item($r, 1).Value = 45828
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 85
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.1
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Neutral"
$ws.Cells.Item($r, 12).Value = 7
$ws.Cells.Item($r, 13).Value = 0.65
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 29.98
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 0.2
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 50
$ws.Cells.Item($r, 20).Value = 12

$r = 291
$ws.Cells.Item($r, 1).Value = 45828
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Small"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 85
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.2
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Neutral"
$ws.Cells.Item($r, 12).Value = 7
$ws.Cells.Item($r, 13).Value = 0.65
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 29.98
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 0.2
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 50
$ws.Cells.Item($r, 20).Value = 12

$r = 292
$ws.Cells.Item($r, 1).Value = 45828
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 85
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.2
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Neutral"
$ws.Cells.Item($r, 12).Value = 7
$ws.Cells.Item($r, 13).Value = 0.65
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 29.98
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 0.2
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 50
$ws.Cells.Item($r, 20).Value = 12

$r = 293
$ws.Cells.Item($r, 1).Value = 45828
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 85
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.5
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 7
$ws.Cells.Item($r, 13).Value = 0.65
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 29.98
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 0.2
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 50
$ws.Cells.Item($r, 20).Value = 12

$r = 294
$ws.Cells.Item($r, 1).Value = 45828
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Large"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 85
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 4
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 7
$ws.Cells.Item($r, 13).Value = 0.65
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 29.98
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 0.2
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 50
$ws.Cells.Item($r, 20).Value = 12

$r = 295
$ws.Cells.Item($r, 1).Value = 45828
$ws.Cells.Item($r, 2).Value = "Tree"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 85
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 1
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 7
$ws.Cells.Item($r, 13).Value = 0.65
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 29.98
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 0.2
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 50
$ws.Cells.Item($r, 20).Value = 12

$r = 296
$ws.Cells.Item($r, 1).Value = 45829
$ws.Cells.Item($r, 2).Value = "Flowering"
$ws.Cells.Item($r, 3).Value = "Large"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 87
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "Yes"
$ws.Cells.Item($r, 10).Value = 2
$ws.Cells.Item($r, 11).Value = "Neutral"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.58
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 30.02
$ws.Cells.Item($r, 16).Value = 15
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 45
$ws.Cells.Item($r, 20).Value = 9

$r = 297
$ws.Cells.Item($r, 1).Value = 45829
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 87
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.1
$ws.Cells.Item($r, 9).Value = "Yes"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.58
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 30.02
$ws.Cells.Item($r, 16).Value = 15
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 45
$ws.Cells.Item($r, 20).Value = 9

$r = 298
$ws.Cells.Item($r, 1).Value = 45829
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Small"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 87
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.2
$ws.Cells.Item($r, 9).Value = "Yes"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Neutral"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.58
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 30.02
$ws.Cells.Item($r, 16).Value = 15
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 45
$ws.Cells.Item($r, 20).Value = 9

$r = 299
$ws.Cells.Item($r, 1).Value = 45829
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 87
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.4
$ws.Cells.Item($r, 9).Value = "Yes"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.58
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 30.02
$ws.Cells.Item($r, 16).Value = 15
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 45
$ws.Cells.Item($r, 20).Value = 9

$r = 300
$ws.Cells.Item($r, 1).Value = 45829
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 87
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.5
$ws.Cells.Item($r, 9).Value = "Yes"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.58
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 30.02
$ws.Cells.Item($r, 16).Value = 15
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 45
$ws.Cells.Item($r, 20).Value = 9

$r = 301
$ws.Cells.Item($r, 1).Value = 45829
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Large"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 87
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "Yes"
$ws.Cells.Item($r, 10).Value = 4
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.58
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 30.02
$ws.Cells.Item($r, 16).Value = 15
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 45
$ws.Cells.Item($r, 20).Value = 9

$r = 302
$ws.Cells.Item($r, 1).Value = 45829
$ws.Cells.Item($r, 2).Value = "Tree"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 70
$ws.Cells.Item($r, 5).Value = 87
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "Yes"
$ws.Cells.Item($r, 10).Value = 1
$ws.Cells.Item($r, 11).Value = "Neutral"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.58
$ws.Cells.Item($r, 14).Value = 70
$ws.Cells.Item($r, 15).Value = 30.02
$ws.Cells.Item($r, 16).Value = 15
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 45
$ws.Cells.Item($r, 20).Value = 9

$r = 303
$ws.Cells.Item($r, 1).Value = 45830
$ws.Cells.Item($r, 2).Value = "Flowering"
$ws.Cells.Item($r, 3).Value = "Large"
$ws.Cells.Item($r, 4).Value = 72
$ws.Cells.Item($r, 5).Value = 92
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.1
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 2
$ws.Cells.Item($r, 11).Value = "Neutral"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.56
$ws.Cells.Item($r, 14).Value = 74
$ws.Cells.Item($r, 15).Value = 30.13
$ws.Cells.Item($r, 16).Value = 12
$ws.Cells.Item($r, 17).Value = 0.11
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 53
$ws.Cells.Item($r, 20).Value = 10

$r = 304
$ws.Cells.Item($r, 1).Value = 45830
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 72
$ws.Cells.Item($r, 5).Value = 92
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.4
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Neutral"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.56
$ws.Cells.Item($r, 14).Value = 74
$ws.Cells.Item($r, 15).Value = 30.13
$ws.Cells.Item($r, 16).Value = 12
$ws.Cells.Item($r, 17).Value = 0.11
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 53
$ws.Cells.Item($r, 20).Value = 10

$r = 305
$ws.Cells.Item($r, 1).Value = 45830
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Small"
$ws.Cells.Item($r, 4).Value = 72
$ws.Cells.Item($r, 5).Value = 92
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.56
$ws.Cells.Item($r, 14).Value = 74
$ws.Cells.Item($r, 15).Value = 30.13
$ws.Cells.Item($r, 16).Value = 12
$ws.Cells.Item($r, 17).Value = 0.11
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 53
$ws.Cells.Item($r, 20).Value = 10

$r = 306
$ws.Cells.Item($r, 1).Value = 45830
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 72
$ws.Cells.Item($r, 5).Value = 92
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.2
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Dark"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.56
$ws.Cells.Item($r, 14).Value = 74
$ws.Cells.Item($r, 15).Value = 30.13
$ws.Cells.Item($r, 16).Value = 12
$ws.Cells.Item($r, 17).Value = 0.11
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 53
$ws.Cells.Item($r, 20).Value = 10

$r = 307
$ws.Cells.Item($r, 1).Value = 45830
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 72
$ws.Cells.Item($r, 5).Value = 92
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = "Neutral"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.56
$ws.Cells.Item($r, 14).Value = 74
$ws.Cells.Item($r, 15).Value = 30.13
$ws.Cells.Item($r, 16).Value = 12
$ws.Cells.Item($r, 17).Value = 0.11
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 53
$ws.Cells.Item($r, 20).Value = 10

$r = 308
$ws.Cells.Item($r, 1).Value = 45830
$ws.Cells.Item($r, 2).Value = "Nonflowering"
$ws.Cells.Item($r, 3).Value = "Large"
$ws.Cells.Item($r, 4).Value = 72
$ws.Cells.Item($r, 5).Value = 92
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0.5
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 4
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.56
$ws.Cells.Item($r, 14).Value = 74
$ws.Cells.Item($r, 15).Value = 30.13
$ws.Cells.Item($r, 16).Value = 12
$ws.Cells.Item($r, 17).Value = 0.11
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 53
$ws.Cells.Item($r, 20).Value = 10

$r = 309
$ws.Cells.Item($r, 1).Value = 45830
$ws.Cells.Item($r, 2).Value = "Tree"
$ws.Cells.Item($r, 3).Value = "Medium"
$ws.Cells.Item($r, 4).Value = 72
$ws.Cells.Item($r, 5).Value = 92
$ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Value = "No"
$ws.Cells.Item($r, 10).Value = 1
$ws.Cells.Item($r, 11).Value = "Bright"
$ws.Cells.Item($r, 12).Value = 9
$ws.Cells.Item($r, 13).Value = 0.56
$ws.Cells.Item($r, 14).Value = 74
$ws.Cells.Item($r, 15).Value = 30.13
$ws.Cells.Item($r, 16).Value = 12
$ws.Cells.Item($r, 17).Value = 0.11
$ws.Cells.Item($r, 18).Value = 9.9
$ws.Cells.Item($r, 19).Value = 53
$ws.Cells.Item($r, 20).Value = 10

# Match the recorded selection state after the edit
$ws.Range("U2").Select()